$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "CanClone" (column M) for the City scene row (row 6) was mis-configured
# as 0; fix it to 1.
$ws.Range("M6").Value = 1

# Restore the cursor/selection to where the editor last left it.
$ws.Range("N12").Select() | Out-Null
